# Halite Helper update of 03/12
# Updates the sample Entity1/Entity2 coordinates, adds a win/loss confidence
# interval block (rows 7-14) and a score block (rows 17-19), removing the
# old scratch row 6 content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing sample data (row 2 / row 3) ---
$ws.Range("B2").Value = 152.99
$ws.Range("C2").Value = 87.06
$ws.Range("B3").Value = 154
$ws.Range("C3").Value = 84

# --- Remove the old scratch row 6 ---
$ws.Range("D6").ClearContents()
$ws.Range("E6").ClearContents()

# --- Labels first, in the order they were authored (fixes shared-string order) ---
$ws.Range("A7").Value = "samples"
$ws.Range("A8").Value = "won"
$ws.Range("A9").Value = "lost"
$ws.Range("A13").Value = "lower bound"
$ws.Range("A14").Value = "upper bound"
$ws.Range("A12").Value = "average"
$ws.Range("A10").Value = "stddev"
$ws.Range("A19").Value = "score"
$ws.Range("A18").Value = "std"
$ws.Range("A17").Value = "average"

# --- New block: win-rate confidence interval (rows 7-14) ---
$ws.Range("B7").Value = 30
$ws.Range("B8").Value = 18
$ws.Range("B9").Formula = "=B7-B8"
$ws.Range("B10").Value = 1.96
$ws.Range("B12").Formula = "=B8/B7"
$ws.Range("B13").Formula = "=B12-B10*SQRT(B12*(1-B12)/B7)"
$ws.Range("B14").Formula = "=B12+B10*SQRT(B12*(1-B12)/B7)"

# --- New block: score (rows 17-19) ---
$ws.Range("B17").Value = 47.56
$ws.Range("B18").Value = 0.5
$ws.Range("B19").Formula = "=B17-3*B18"

# --- Column widths (best-effort match of the autosized columns) ---
$ws.Columns.Item(1).ColumnWidth = 11.592447916666666
$ws.Columns.Item(5).ColumnWidth = 11.166666666666666

# --- Selection moves to C20 after the edits ---
$null = $ws.Range("C20").Select()
